# Scheduled market-data refresh: update computed price/profit columns (H:N)
# across the Leve-profit worksheets, per the latest Universalis price pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Range("H18").Value = 4524.25
$ws.Range("I18").Value = 4524.25
$ws.Range("K18").Value = 4524.25
$ws.Range("M18").Value = -4240.25

# Row 69
$ws.Range("H69").Value = 3814.5186
$ws.Range("I69").Value = 1517.8182
$ws.Range("J69").Value = 5393.5
$ws.Range("K69").Value = 4553.4546
$ws.Range("L69").Value = 16180.5
$ws.Range("M69").Value = -3679.4546
$ws.Range("N69").Value = -17928.5

# Row 72
$ws.Range("H72").Value = 3814.5186
$ws.Range("I72").Value = 1517.8182
$ws.Range("J72").Value = 5393.5
$ws.Range("K72").Value = 13660.3638
$ws.Range("L72").Value = 48541.5
$ws.Range("M72").Value = -9292.363799999999
$ws.Range("N72").Value = -57277.5

# Row 112
$ws.Range("H112").Value = 68425.47
$ws.Range("J112").Value = 68425.47
$ws.Range("L112").Value = 205276.41
$ws.Range("N112").Value = -207492.41

# Row 116
$ws.Range("H116").Value = 3086.875
$ws.Range("I116").Value = 2050
$ws.Range("J116").Value = 4123.75
$ws.Range("K116").Value = 2050
$ws.Range("L116").Value = 4123.75
$ws.Range("M116").Value = 1392
$ws.Range("N116").Value = -11007.75

# Row 138
$ws.Range("H138").Value = 5029.203
$ws.Range("I138").Value = 7999.227
$ws.Range("K138").Value = 23997.681
$ws.Range("M138").Value = -18857.681

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 190495.5
$ws.Range("I32").Value = 219624.88
$ws.Range("K32").Value = 219624.88
$ws.Range("M32").Value = -219337.88

# Row 45
$ws.Range("H45").Value = 93616.63
$ws.Range("I45").Value = 93616.63
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 93616.63
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -93239.63
$ws.Range("N45").Value = $null

# Row 64
$ws.Range("H64").Value = 50000
$ws.Range("J64").Value = 50000
$ws.Range("L64").Value = 50000
$ws.Range("N64").Value = -50496

# Row 67
$ws.Range("H67").Value = 50000
$ws.Range("J67").Value = 50000
$ws.Range("L67").Value = 50000
$ws.Range("N67").Value = -51716

# Row 102
$ws.Range("H102").Value = 3169.5715
$ws.Range("I102").Value = 2396.25
$ws.Range("J102").Value = 4200.6665
$ws.Range("K102").Value = 2396.25
$ws.Range("L102").Value = 4200.6665
$ws.Range("M102").Value = -774.25
$ws.Range("N102").Value = -7444.6665

# Row 122
$ws.Range("H122").Value = 3170.5715
$ws.Range("I122").Value = 2446.4
$ws.Range("K122").Value = 7339.200000000001
$ws.Range("M122").Value = -4889.200000000001

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1061.8889
$ws.Range("I20").Value = 917.55
$ws.Range("J20").Value = 1474.2858
$ws.Range("K20").Value = 917.55
$ws.Range("L20").Value = 1474.2858
$ws.Range("M20").Value = -670.55
$ws.Range("N20").Value = -1968.2858

# Row 86
$ws.Range("H86").Value = 3678.9285
$ws.Range("I86").Value = 2202.75
$ws.Range("J86").Value = 5647.1665
$ws.Range("K86").Value = 2202.75
$ws.Range("L86").Value = 5647.1665
$ws.Range("M86").Value = -1079.75
$ws.Range("N86").Value = -7893.1665

# Row 89
$ws.Range("H89").Value = 3678.9285
$ws.Range("I89").Value = 2202.75
$ws.Range("J89").Value = 5647.1665
$ws.Range("K89").Value = 11013.75
$ws.Range("L89").Value = 28235.8325
$ws.Range("M89").Value = -5397.75
$ws.Range("N89").Value = -39467.8325

# Row 99
$ws.Range("H99").Value = 5770.625
$ws.Range("I99").Value = 11089.556
$ws.Range("K99").Value = 11089.556
$ws.Range("M99").Value = -9591.556

# Row 134
$ws.Range("H134").Value = 18369872
$ws.Range("I134").Value = 1843.8889
$ws.Range("J134").Value = 69235176
$ws.Range("K134").Value = 5531.6667
$ws.Range("L134").Value = 207705528
$ws.Range("M134").Value = -2996.6667
$ws.Range("N134").Value = -207710598

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2926.0317
$ws.Range("I31").Value = 1905.9
$ws.Range("J31").Value = 3118.5095
$ws.Range("K31").Value = 1905.9
$ws.Range("L31").Value = 3118.5095
$ws.Range("M31").Value = -1610.9
$ws.Range("N31").Value = -3708.5095

# Row 34
$ws.Range("H34").Value = 2926.0317
$ws.Range("I34").Value = 1905.9
$ws.Range("J34").Value = 3118.5095
$ws.Range("K34").Value = 1905.9
$ws.Range("L34").Value = 3118.5095
$ws.Range("M34").Value = -1703.9
$ws.Range("N34").Value = -3522.5095

# Row 58
$ws.Range("H58").Value = 3176.9583
$ws.Range("I58").Value = 2713.7273
$ws.Range("J58").Value = 3568.923
$ws.Range("K58").Value = 2713.7273
$ws.Range("L58").Value = 3568.923
$ws.Range("M58").Value = -2510.7273
$ws.Range("N58").Value = -3974.923

# Row 62
$ws.Range("H62").Value = 11280.363
$ws.Range("I62").Value = 16183.571
$ws.Range("J62").Value = 2699.75
$ws.Range("K62").Value = 16183.571
$ws.Range("L62").Value = 2699.75
$ws.Range("M62").Value = -15559.571
$ws.Range("N62").Value = -3947.75

# Row 65
$ws.Range("H65").Value = 11280.363
$ws.Range("I65").Value = 16183.571
$ws.Range("J65").Value = 2699.75
$ws.Range("K65").Value = 80917.855
$ws.Range("L65").Value = 13498.75
$ws.Range("M65").Value = -77797.855
$ws.Range("N65").Value = -19738.75

# Row 107
$ws.Range("H107").Value = 2061.2727
$ws.Range("I107").Value = 1609.5625
$ws.Range("J107").Value = 3265.8333
$ws.Range("K107").Value = 1609.5625
$ws.Range("L107").Value = 3265.8333
$ws.Range("M107").Value = 310.4375
$ws.Range("N107").Value = -7105.8333

# Row 136
$ws.Range("H136").Value = 3176.9583
$ws.Range("I136").Value = 2713.7273
$ws.Range("J136").Value = 3568.923
$ws.Range("K136").Value = 8141.1819
$ws.Range("L136").Value = 10706.769
$ws.Range("M136").Value = -5591.1819
$ws.Range("N136").Value = -15806.769

$ws = $wb.Worksheets.Item("CUL")
# Row 23
$ws.Range("H23").Value = 333399.66
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").Value = $null

# Row 33
$ws.Range("H33").Value = 57269.168
$ws.Range("I33").Value = 1116.6666
$ws.Range("J33").Value = 68499.664
$ws.Range("K33").Value = 6699.9996
$ws.Range("L33").Value = 410997.9840000001
$ws.Range("M33").Value = -6416.9996
$ws.Range("N33").Value = -411563.9840000001

# Row 39
$ws.Range("H39").Value = 6633.25
$ws.Range("I39").Value = 855
$ws.Range("J39").Value = 7966.6924
$ws.Range("K39").Value = 2565
$ws.Range("L39").Value = 23900.0772
$ws.Range("M39").Value = -2271
$ws.Range("N39").Value = -24488.0772

# Row 80
$ws.Range("H80").Value = 1923.5
$ws.Range("J80").Value = 2002.5
$ws.Range("L80").Value = 6007.5
$ws.Range("N80").Value = -7879.5

# Row 81
$ws.Range("H81").Value = 41671412
$ws.Range("I81").Value = 166667420
$ws.Range("K81").Value = 500002260
$ws.Range("M81").Value = -500001137

# Row 83
$ws.Range("H83").Value = 1923.5
$ws.Range("J83").Value = 2002.5
$ws.Range("L83").Value = 18022.5
$ws.Range("N83").Value = -27382.5

# Row 84
$ws.Range("H84").Value = 41671412
$ws.Range("I84").Value = 166667420
$ws.Range("K84").Value = 1500006780
$ws.Range("M84").Value = -1500001164

# Row 87
$ws.Range("H87").Value = 13508
$ws.Range("I87").Value = 3117.25
$ws.Range("K87").Value = 9351.75
$ws.Range("M87").Value = -8103.75

# Row 88
$ws.Range("H88").Value = 9118.462
$ws.Range("I88").Value = 4875
$ws.Range("K88").Value = 14625
$ws.Range("M88").Value = -14197

# Row 90
$ws.Range("H90").Value = 13508
$ws.Range("I90").Value = 3117.25
$ws.Range("K90").Value = 28055.25
$ws.Range("M90").Value = -21815.25

# Row 91
$ws.Range("H91").Value = 9118.462
$ws.Range("I91").Value = 4875
$ws.Range("K91").Value = 14625
$ws.Range("M91").Value = -13143

# Row 94
$ws.Range("H94").Value = 13633
$ws.Range("I94").Value = 9000
$ws.Range("J94").Value = 15949.5
$ws.Range("K94").Value = 27000
$ws.Range("L94").Value = 47848.5
$ws.Range("M94").Value = -26324
$ws.Range("N94").Value = -49200.5

# Row 97
$ws.Range("H97").Value = 185.16667
$ws.Range("I97").Value = 244
$ws.Range("K97").Value = 732
$ws.Range("M97").Value = -236

# Row 100
$ws.Range("H100").Value = 4666.6665
$ws.Range("J100").Value = 4666.6665
$ws.Range("L100").Value = 13999.9995
$ws.Range("N100").Value = -15621.9995

# Row 103
$ws.Range("H103").Value = 1258623.5
$ws.Range("I103").Value = 3334333
$ws.Range("J103").Value = 13197.8
$ws.Range("K103").Value = 10002999
$ws.Range("L103").Value = 39593.39999999999
$ws.Range("M103").Value = -10002120
$ws.Range("N103").Value = -41351.39999999999

# Row 106
$ws.Range("H106").Value = 12496.333
$ws.Range("J106").Value = 16244.5
$ws.Range("L106").Value = 48733.5
$ws.Range("N106").Value = -50625.5

# Row 112
$ws.Range("H112").Value = 10624.875
$ws.Range("I112").Value = 8500
$ws.Range("J112").Value = 12749.75
$ws.Range("K112").Value = 25500
$ws.Range("L112").Value = 38249.25
$ws.Range("M112").Value = -24392
$ws.Range("N112").Value = -40465.25

# Row 134
$ws.Range("H134").Value = 17859892
$ws.Range("I134").Value = 20835708
$ws.Range("K134").Value = 62507124
$ws.Range("M134").Value = -62502054

# Row 139
$ws.Range("H139").Value = 4632551.5
$ws.Range("J139").Value = 3734.0588
$ws.Range("L139").Value = 11202.1764
$ws.Range("N139").Value = -21482.1764

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 15389707
$ws.Range("I61").Value = 20005120
$ws.Range("J61").Value = 4998.3335
$ws.Range("K61").Value = 20005120
$ws.Range("L61").Value = 4998.3335
$ws.Range("M61").Value = -20004918
$ws.Range("N61").Value = -5402.3335

# Row 113
$ws.Range("H113").Value = 15389707
$ws.Range("I113").Value = 20005120
$ws.Range("J113").Value = 4998.3335
$ws.Range("K113").Value = 20005120
$ws.Range("L113").Value = 4998.3335
$ws.Range("M113").Value = -20002950
$ws.Range("N113").Value = -9338.333500000001

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 517.7778
$ws.Range("I113").Value = 501.4
$ws.Range("K113").Value = 1504.2
$ws.Range("M113").Value = 665.8000000000002

# Row 135
$ws.Range("H135").Value = 85039.2
$ws.Range("J135").Value = 85039.2
$ws.Range("L135").Value = 85039.2
$ws.Range("N135").Value = -95179.2

# Row 136
$ws.Range("H136").Value = 50397.047
$ws.Range("I136").Value = 126493.625
$ws.Range("J136").Value = 3568.3845
$ws.Range("K136").Value = 379480.875
$ws.Range("L136").Value = 10705.1535
$ws.Range("M136").Value = -376930.875
$ws.Range("N136").Value = -15805.1535
